{"js": "// Target edit (per the commit's XML diff):\n//   - Add a new blank paragraph at the very top of the body.\n//   - Add a new paragraph containing \"fkshfkhks\" right after it (i.e. right\n//     before the existing \"s,afkfaf\" paragraph).\n//   - Move the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd pair) from the\n//     end of the \"s,afkfaf\" paragraph (after its run) to the start of that\n//     same paragraph (before its run).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The original (and, at this point, only) paragraph: \"s,afkfaf\".\nconst originalFirst = paragraphs.items[0];\n\n// Insert the new \"fkshfkhks\" paragraph immediately before the original one.\nconst middlePara = originalFirst.insertParagraph(\"fkshfkhks\", Word.InsertLocation.before);\n\n// Insert a new blank paragraph before that \"fkshfkhks\" paragraph, so it\n// becomes the first paragraph in the document.\nmiddlePara.insertParagraph(\"\", Word.InsertLocation.before);\nawait context.sync();\n\n// Relocate the \"_GoBack\" bookmark from the end of the original paragraph to\n// its start: delete the old bookmark, then insert a fresh one (same name)\n// at the very beginning of that paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nconst startRange = originalFirst.getRange(Word.RangeLocation.start);\nstartRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Target edit (per the commit's XML diff):\n#   - Add a new blank paragraph at the very top of the document.\n#   - Add a new paragraph containing \"fkshfkhks\" right after it (i.e. right\n#     before the existing \"s,afkfaf\" paragraph).\n#   - Move the \"_GoBack\" bookmark (bookmarkStart/bookmarkEnd pair) from the\n#     end of the \"s,afkfaf\" paragraph (after its run) to the start of that\n#     same paragraph (before its run).\n\n$d = $word.ActiveDocument\n\n# The original (and, at this point, only) paragraph: \"s,afkfaf\".\n$originalFirst = $d.Paragraphs.Item(1)\n$r1 = $originalFirst.Range\n\n# Insert the new \"fkshfkhks\" paragraph immediately before the original one.\n$r1.InsertParagraphBefore()\n$r1.InsertBefore(\"fkshfkhks\")\n\n# Insert a new blank paragraph before that \"fkshfkhks\" paragraph, so it\n# becomes the first paragraph in the document.\n$r1.InsertParagraphBefore()\n\n# Relocate the \"_GoBack\" bookmark from the end of the original paragraph to\n# its start: delete the old bookmark, then add a fresh one (same name) at\n# the very beginning of that paragraph (now the 3rd paragraph).\n$bm = $d.Bookmarks.Item(\"_GoBack\")\n$bm.Delete()\n\n$targetPara = $d.Paragraphs.Item(3)\n$startRange = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)\n$d.Bookmarks.Add(\"_GoBack\", $startRange)\n"}
